{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\n// Find the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph.\nlet target = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text.indexOf(\"Docente(s) Respons\u00e1vel(eis)\") !== -1) {\n    target = paras.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find 'Docente(s) Respons\u00e1vel(eis)' paragraph\");\n}\n\n// Insert a new paragraph right after it, with the ListBullet style.\nconst newPara = target.insertParagraph(\"5701460 - Antonio Iacono\", \"After\");\nnewPara.style = \"List Bullet\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"Docente(s) Respons\u00e1vel(eis)*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find 'Docente(s) Respons\u00e1vel(eis)' paragraph\"\n}\n\n$target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n$newPara.Range.Text = \"5701460 - Antonio Iacono\"\n$newPara.Style = \"List Bullet\"\n"}
